$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Fleet Code" as new column A (shifts everything right) ---
[void]$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "Fleet Code"

# --- Insert "Agent Name" as new column C (after Vehicle Number, before Fitness Number) ---
[void]$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "Agent Name"

# --- Append the new trailing columns (M..S) ---
$ws.Range("M1").Value = "Engine No."
$ws.Range("N1").Value = "Chassis No"
$ws.Range("O1").Value = "Manufacture Year"
$ws.Range("P1").Value = "Type Of Body"
$ws.Range("Q1").Value = "Type Of Fuel"
$ws.Range("R1").Value = "Seating Capacity(including Driver)"
$ws.Range("S1").Value = "Cubic Capacity"

# --- Column widths (B, D, E, F, G keep the exact widths they already had
#     before the column inserts shifted them over, so they are left alone) ---
$ws.Columns.Item(1).ColumnWidth  = 12.833333333333334   # A Fleet Code
$ws.Columns.Item(3).ColumnWidth  = 26.833333333333332   # C Agent Name
$ws.Columns.Item(8).ColumnWidth  = 17.5                 # H Pay Date
$ws.Columns.Item(9).ColumnWidth  = 21.0                 # I Pay Bank
$ws.Columns.Item(10).ColumnWidth = 19.666666666666668   # J Pay Branch
$ws.Columns.Item(11).ColumnWidth = 18.0                 # K Valid From
$ws.Columns.Item(12).ColumnWidth = 20.166666666666668   # L Valid Till
$ws.Columns.Item(13).ColumnWidth = 19.833333333333332   # M Engine No.
$ws.Columns.Item(14).ColumnWidth = 23.666666666666668   # N Chassis No
$ws.Columns.Item(15).ColumnWidth = 28.166666666666668   # O Manufacture Year
$ws.Columns.Item(16).ColumnWidth = 21.666666666666668   # P Type Of Body
$ws.Columns.Item(17).ColumnWidth = 20.666666666666668   # Q Type Of Fuel
$ws.Columns.Item(18).ColumnWidth = 36.166666666666664   # R Seating Capacity
$ws.Columns.Item(19).ColumnWidth = 18.666666666666668   # S Cubic Capacity

# --- Row 1 height ---
$ws.Rows.Item(1).RowHeight = 15

# --- Selection / view state ---
[void]$ws.Range("C1:C1048576").Select()

Write-Host "Applied Fleet Code / Agent Name / trailing vehicle-detail columns."
